$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H39").Value = 888.5
$ws.Range("I39").Value = 100.666664
$ws.Range("J39").Value = 1901.4286
$ws.Range("K39").Value = 301.999992
$ws.Range("L39").Value = 5704.2858
$ws.Range("M39").Value = -5.99999200000002
$ws.Range("N39").Value = -6296.2858
$ws.Range("H93").Value = 85450.5
$ws.Range("J93").Value = 85450.5
$ws.Range("L93").Value = 85450.5
$ws.Range("N93").Value = -90442.5
$ws.Range("H116").Value = 14726.223
$ws.Range("J116").Value = 4067
$ws.Range("L116").Value = 4067
$ws.Range("N116").Value = -10951
$ws.Range("H121").Value = 1729
$ws.Range("J121").Value = 2493.5
$ws.Range("L121").Value = 7480.5
$ws.Range("N121").Value = -10974.5
$ws.Range("H130").Value = 39593.332
$ws.Range("J130").Value = 39593.332
$ws.Range("L130").Value = 39593.332
$ws.Range("N130").Value = -49633.332
$ws.Range("H132").Value = 735.1316
$ws.Range("I132").Value = 619.1818
$ws.Range("K132").Value = 1857.5454
$ws.Range("M132").Value = 672.4546
$ws.Range("H135").Value = 339.70587
$ws.Range("H137").Value = 2129.6
$ws.Range("I137").Value = 1910.1428
$ws.Range("K137").Value = 5730.428400000001
$ws.Range("M137").Value = -3180.428400000001
$ws.Range("H138").Value = 1808.9836
$ws.Range("I138").Value = 1342.125
$ws.Range("K138").Value = 4026.375
$ws.Range("M138").Value = 1113.625
$ws.Range("H141").Value = 4005063.2
$ws.Range("I141").Value = 5601397
$ws.Range("J141").Value = 14228
$ws.Range("K141").Value = 16804191
$ws.Range("L141").Value = 42684
$ws.Range("M141").Value = -16799011
$ws.Range("N141").Value = -53044

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5556.4033
$ws.Range("I32").Value = 3696.8333
$ws.Range("J32").Value = 15474.111
$ws.Range("K32").Value = 3696.8333
$ws.Range("L32").Value = 15474.111
$ws.Range("M32").Value = -3409.8333
$ws.Range("N32").Value = -16048.111
$ws.Range("H38").Value = 0
$ws.Range("I38").Value = 0
$ws.Range("K38").Value = 0
$ws.Range("M38").ClearContents()
$ws.Range("H41").Value = 3022
$ws.Range("I41").Value = 3022
$ws.Range("J41").Value = 0
$ws.Range("K41").Value = 3022
$ws.Range("L41").Value = 0
$ws.Range("M41").Value = -2608
$ws.Range("N41").ClearContents()
$ws.Range("H43").Value = 38900
$ws.Range("J43").Value = 38900
$ws.Range("L43").Value = 38900
$ws.Range("N43").Value = -39526
$ws.Range("H45").Value = 4501301
$ws.Range("J45").Value = 1679.2727
$ws.Range("L45").Value = 1679.2727
$ws.Range("N45").Value = -2433.2727
$ws.Range("H74").Value = 1236.4546
$ws.Range("I74").Value = 750.2143
$ws.Range("K74").Value = 750.2143
$ws.Range("M74").Value = 123.7857
$ws.Range("H77").Value = 1236.4546
$ws.Range("I77").Value = 750.2143
$ws.Range("K77").Value = 3751.0715
$ws.Range("M77").Value = 616.9285
$ws.Range("H97").Value = 974.4286
$ws.Range("I97").Value = 943.2
$ws.Range("J97").Value = 1599
$ws.Range("K97").Value = 943.2
$ws.Range("L97").Value = 1599
$ws.Range("M97").Value = -447.2
$ws.Range("N97").Value = -2591
$ws.Range("H110").Value = 1659.5714
$ws.Range("I110").Value = 269.5
$ws.Range("J110").Value = 10000
$ws.Range("K110").Value = 269.5
$ws.Range("L110").Value = 10000
$ws.Range("M110").Value = 1775.5
$ws.Range("N110").Value = -14090
$ws.Range("H132").Value = 1493.0526
$ws.Range("I132").Value = 1174.1666
$ws.Range("K132").Value = 3522.4998
$ws.Range("M132").Value = -992.4998000000001
$ws.Range("H137").Value = 29999
$ws.Range("J137").Value = 29999
$ws.Range("L137").Value = 29999
$ws.Range("N137").Value = -40199

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 4666
$ws.Range("I20").Value = 4499
$ws.Range("K20").Value = 4499
$ws.Range("M20").Value = -4252
$ws.Range("H117").Value = 56000
$ws.Range("J117").Value = 56000
$ws.Range("L117").Value = 56000
$ws.Range("N117").Value = -65178
$ws.Range("H128").Value = 4833.3335
$ws.Range("I128").Value = 4833.3335
$ws.Range("K128").Value = 14500.0005
$ws.Range("M128").Value = -12010.0005

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H19").Value = 860
$ws.Range("I19").Value = 860
$ws.Range("K19").Value = 860
$ws.Range("M19").Value = -690
$ws.Range("H24").Value = 860
$ws.Range("I24").Value = 860
$ws.Range("K24").Value = 860
$ws.Range("M24").Value = -690
$ws.Range("H53").Value = 52583
$ws.Range("J53").Value = 52583
$ws.Range("L53").Value = 52583
$ws.Range("N53").Value = -53797
$ws.Range("H134").Value = 3352.75
$ws.Range("J134").Value = 4608.6665
$ws.Range("L134").Value = 13825.9995
$ws.Range("N134").Value = -18895.9995
$ws.Range("H135").Value = 39850.2
$ws.Range("J135").Value = 39850.2
$ws.Range("L135").Value = 39850.2
$ws.Range("N135").Value = -49990.2

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 11269.539
$ws.Range("J131").Value = 11553.446
$ws.Range("L131").Value = 34660.338
$ws.Range("N131").Value = -44740.338

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 140
$ws.Range("I2").Value = 83
$ws.Range("J2").Value = 182.75
$ws.Range("K2").Value = 83
$ws.Range("L2").Value = 182.75
$ws.Range("M2").Value = 30
$ws.Range("N2").Value = -408.75
$ws.Range("H20").Value = 1578001.2
$ws.Range("J20").Value = 12000
$ws.Range("L20").Value = 12000
$ws.Range("N20").Value = -12490
$ws.Range("H24").Value = 3343333.2
$ws.Range("J24").Value = 15000
$ws.Range("L24").Value = 15000
$ws.Range("N24").Value = -15346
$ws.Range("H48").Value = 19976.666
$ws.Range("J48").Value = 19976.666
$ws.Range("L48").Value = 19976.666
$ws.Range("N48").Value = -20946.666
$ws.Range("H70").Value = 0
$ws.Range("I70").Value = 0
$ws.Range("J70").Value = 0
$ws.Range("K70").Value = 0
$ws.Range("L70").Value = 0
$ws.Range("M70").ClearContents()
$ws.Range("N70").ClearContents()
$ws.Range("H73").Value = 0
$ws.Range("I73").Value = 0
$ws.Range("J73").Value = 0
$ws.Range("K73").Value = 0
$ws.Range("L73").Value = 0
$ws.Range("M73").ClearContents()
$ws.Range("N73").ClearContents()
$ws.Range("H92").Value = 21879.8
$ws.Range("J92").Value = 21879.8
$ws.Range("L92").Value = 21879.8
$ws.Range("N92").Value = -25623.8
$ws.Range("H97").Value = 607.89655
$ws.Range("I97").Value = 620.7037
$ws.Range("K97").Value = 620.7037
$ws.Range("M97").Value = -124.7037
$ws.Range("H101").Value = 0
$ws.Range("J101").Value = 0
$ws.Range("L101").Value = 0
$ws.Range("N101").ClearContents()
$ws.Range("H102").Value = 2609.4
$ws.Range("I102").Value = 2326.4666
$ws.Range("J102").Value = 3458.2
$ws.Range("K102").Value = 2326.4666
$ws.Range("L102").Value = 3458.2
$ws.Range("M102").Value = -704.4666000000002
$ws.Range("N102").Value = -6702.2
$ws.Range("H113").Value = 1433.125
$ws.Range("I113").Value = 1147.6666
$ws.Range("K113").Value = 1147.6666
$ws.Range("M113").Value = 1022.3334
$ws.Range("H132").Value = 2139248
$ws.Range("I132").Value = 3206796
$ws.Range("J132").Value = 4151.6665
$ws.Range("K132").Value = 9620388
$ws.Range("L132").Value = 12454.9995
$ws.Range("M132").Value = -9617858
$ws.Range("N132").Value = -17514.9995

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H32").Value = 9870
$ws.Range("I32").Value = 0
$ws.Range("J32").Value = 9870
$ws.Range("K32").Value = 0
$ws.Range("L32").Value = 9870
$ws.Range("M32").ClearContents()
$ws.Range("N32").Value = -10504
$ws.Range("H46").Value = 1652.0625
$ws.Range("I46").Value = 985.2
$ws.Range("K46").Value = 985.2
$ws.Range("M46").Value = -797.2
$ws.Range("H55").Value = 558.93335
$ws.Range("I55").Value = 564.8889
$ws.Range("J55").Value = 550
$ws.Range("K55").Value = 564.8889
$ws.Range("L55").Value = 550
$ws.Range("M55").Value = -391.8889
$ws.Range("N55").Value = -896

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H31").Value = 17000
$ws.Range("J31").Value = 17000
$ws.Range("L31").Value = 17000
$ws.Range("N31").Value = -17696
$ws.Range("H132").Value = 988.1594
$ws.Range("I132").Value = 780.8461
$ws.Range("K132").Value = 2342.5383
$ws.Range("M132").Value = 187.4616999999998
